$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks first - we'll rebuild the whole table (including
# the two pre-existing rows, which move further down the sheet) from scratch
# so every hyperlink ends up anchored to the right cell.
$ws.Hyperlinks.Delete()

# Header row (unchanged)
$ws.Range("A1").Value = "Label"
$ws.Range("B1").Value = "URL"
$ws.Range("C1").Value = "Picture"

# Full data set in final row order.
$rows = @(
    @{ Row = 2;  Label = "YouTube";               Url = "https://www.youtube.com";                                                                                                                                   Pic = "Youtube.jpg" },
    @{ Row = 3;  Label = "Free Media Heck Yeah";   Url = "https://fmhy.net";                                                                                                                                          Pic = "FMHY.jpeg" },
    @{ Row = 4;  Label = "Miruro Anima";           Url = "https://www.miruro.tv";                                                                                                                                     Pic = "miruro-tv-icon.jpg" },
    @{ Row = 5;  Label = "Pandora";                Url = "https://www.pandora.com";                                                                                                                                  Pic = "pandora.png" },
    @{ Row = 6;  Label = "Shangri-La Frontier";    Url = "https://www.crunchyroll.com/series/G79H23Z8P/shangri-la-frontier";                                                                                         Pic = "ShangriLaFrontier.png" },
    @{ Row = 7;  Label = "Solo Levelling";         Url = "https://www.crunchyroll.com/series/GDKHZEJ0K/solo-leveling";                                                                                               Pic = "SoloLeveling.webp" },
    @{ Row = 8;  Label = "Demon Slayer";           Url = "https://www.crunchyroll.com/series/GY5P48XEY/demon-slayer-kimetsu-no-yaiba?srsltid=AfmBOoqt8gh7y_FZv2IkKBba9pTdA26VzePJukANMW2uAR_fGGVYSlk7"; Pic = "DemonSlayer.jpg" },
    @{ Row = 9;  Label = "Hulu";                   Url = "https://www.hulu.com/welcome";                                                                                                                             Pic = "hulu.jpg" },
    @{ Row = 10; Label = "Sid the Scientist";      Url = "https://pbskids.org/sid";                                                                                                                                   Pic = "Sid.jpg" },
    @{ Row = 11; Label = "Daniel Tiger";           Url = "https://pbskids.org/daniel/";                                                                                                                               Pic = "DanielTiger.jpg" },
    @{ Row = 12; Label = "PBS Kids Videos";        Url = "https://pbskids.org/videos";                                                                                                                                Pic = "PBS-Kids-Video.png" },
    @{ Row = 13; Label = "PBS Kids Games";         Url = "https://pbskids.org/games";                                                                                                                                 Pic = "PBS_Kids_Games.png" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value = $r.Label
    $ws.Range("B$rowNum").Value = $r.Url
    $ws.Range("C$rowNum").Value = $r.Pic

    $null = $ws.Hyperlinks.Add($ws.Range("B$rowNum"), $r.Url)
    $ws.Range("B$rowNum").Style = "Hyperlink"
}

# Restore the active selection left behind by the edit.
$null = $ws.Range("B14").Select()

# Reposition/resize the workbook window to match the saved view.
$win = $excel.Windows.Item(1)
$win.Left = 14295
$win.Top = 0
$win.Width = 14610
$win.Height = 15585
